$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 / Row 2: shift labels changed from "6am-2pm"/"2pm-10pm" to "A"/"B"
$ws.Range("A1").Value = "A"
$ws.Range("A2").Value = "B"

# New row 3: third shift "C" with Aluddin / Ainal and the same date as rows 1-2
$ws.Range("A3").Value = "C"
$ws.Range("B3").Value = "Aluddin"
$ws.Range("C3").Value = "Ainal"
$ws.Range("D3").Value = 45965

# Copy the date's number formatting (style) from D1 onto D3 so it matches
# the existing date formatting instead of creating a new number format
$ws.Range("D1").Copy()
$null = $ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move / update the active selection to A2
$null = $ws.Range("A2").Select()
